$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: fill in the vendor/price for the existing ASIN (previously spread
# across rows 2 and 3)
$ws.Range("B2").Value = "Amazon.ae"
$ws.Range("C2").Value = "AED183.16"
$ws.Range("D2").Value = "Does not import internationally"

# Row 3: new ASIN row, vendor/price not yet known (mirrors the old row 2) -
# clear the stale vendor/price that used to live here
$ws.Range("A3").Value = "B007177NZU"
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = "Does not import internationally"

# Row 4: new ASIN, Amazon.ae offer
$ws.Range("A4").Value = "B007177NZU"
$ws.Range("B4").Value = "Amazon.ae"
$ws.Range("C4").Value = "AED294.58"
$ws.Range("D4").Value = "Does not import internationally"

# Row 5: new ASIN, Amazon.ae offer
$ws.Range("A5").Value = "B007177NZU"
$ws.Range("B5").Value = "Amazon.ae"
$ws.Range("C5").Value = "AED295.00"
$ws.Range("D5").Value = "Does not import internationally"

# Row 6: new ASIN, Amazon.ae offer
$ws.Range("A6").Value = "B007177NZU"
$ws.Range("B6").Value = "Amazon.ae"
$ws.Range("C6").Value = "AED299.00"
$ws.Range("D6").Value = "Does not import internationally"

# Row 7: new ASIN, Amazon.ae offer (duplicate price of row 6)
$ws.Range("A7").Value = "B007177NZU"
$ws.Range("B7").Value = "Amazon.ae"
$ws.Range("C7").Value = "AED299.00"
$ws.Range("D7").Value = "Does not import internationally"

# Row 8: new ASIN, third-party seller offer shipped from outside the UAE
$ws.Range("A8").Value = "B007177NZU"
$ws.Range("B8").Value = "UNER STORE"
$ws.Range("C8").Value = "AED304.29"
$ws.Range("D8").Value = "Ships from outside the UAE. Learn more"
